# Update automatico via Actualizar 02-14-2021 13-11-49
#
# The "Fecha" (D) column holds a rolling log of check timestamps: each
# refresh inserts a fresh timestamp block at the top (rows 2-15) and the
# previously-newest blocks shift down one block (14 rows) at a time.
#
# New block   (rows 2-15)  -> 2021-02-14 13:11:41 (44241.5497863199)
# Shift block (rows 16-29) -> what used to be in rows 2-15 (44241.52858888889)
# Shift block (rows 30-43) -> what used to be in rows 16-29 (44241.507301875)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:D15").Value = 44241.5497863199
$ws.Range("D16:D29").Value = 44241.52858888889
$ws.Range("D30:D43").Value = 44241.507301875
